# Apply the latest cryptos snapshot update (price + 1h volume change)
# Values in column D are written with a leading apostrophe when they look
# numeric so Excel stores them as text (matching the original inlineStr cells)
# and preserves formatting such as trailing zeros (e.g. "0.4230").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '25.945.31'
$ws.Range("E2").Value = '  -0.95%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.635.28'
$ws.Range("E3").Value = '  -2.46%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.09%  '
# Row 5: BNB
$ws.Range("D5").Value = '''209.22'
$ws.Range("E5").Value = '  -0.98%  '
# Row 6: XRP
$ws.Range("D6").Value = '''0.5166'
$ws.Range("E6").Value = '  -2.09%  '
# Row 7: USDC
$ws.Range("E7").Value = '  +0.11%  '
# Row 8: Cardano
$ws.Range("E8").Value = '  -4.11%  '
# Row 9: Dogecoin
$ws.Range("D9").Value = '''0.06220'
$ws.Range("E9").Value = '  -1.51%  '
# Row 10: Solana
$ws.Range("E10").Value = '  -5.38%  '
# Row 11: TRON
$ws.Range("D11").Value = '''0.07546'
$ws.Range("E11").Value = '  +0.08%  '
# Row 12: WrappedEther
$ws.Range("D12").Value = '1.644.01'
$ws.Range("E12").Value = '  -2.07%  '
# Row 13: Polkadot
$ws.Range("D13").Value = '''4.347'
$ws.Range("E13").Value = '  -2.42%  '
# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '1.862.53'
$ws.Range("E14").Value = '  -2.19%  '
# Row 15: Polygon
$ws.Range("D15").Value = '''0.5403'
$ws.Range("E15").Value = '  -4.51%  '
# Row 16: ShibaInu
$ws.Range("D16").Value = '0.0₅7891'
$ws.Range("E16").Value = '  -1.98%  '
# Row 17: Litecoin
$ws.Range("D17").Value = '''64.39'
$ws.Range("E17").Value = '  -3.27%  '
# Row 18: WrappedBTC
$ws.Range("D18").Value = '25.965.97'
$ws.Range("E18").Value = '  -0.95%  '
# Row 19: Dai
$ws.Range("E19").Value = '  +0.07%  '
# Row 20: Uniswap
$ws.Range("D20").Value = '''4.616'
$ws.Range("E20").Value = '  -4.36%  '
# Row 21: BitcoinCash
$ws.Range("D21").Value = '''184.16'
$ws.Range("E21").Value = '  -2.31%  '
# Row 22: Avalanche
$ws.Range("E22").Value = '  -4.40%  '
# Row 23: Chainlink
$ws.Range("D23").Value = '''6.064'
$ws.Range("E23").Value = '  -2.13%  '
# Row 24: BinanceUSD
$ws.Range("E24").Value = '  +0.13%  '
# Row 25: Monero
$ws.Range("D25").Value = '''145.25'
$ws.Range("E25").Value = '  -2.13%  '
# Row 26: Cosmos
$ws.Range("D26").Value = '''7.302'
$ws.Range("E26").Value = '  -3.91%  '
# Row 27: Stellar
$ws.Range("E27").Value = '  -5.62%  '
# Row 28: EthereumClassic
$ws.Range("D28").Value = '''15.44'
$ws.Range("E28").Value = '  -3.46%  '
# Row 29: Toncoin
$ws.Range("D29").Value = '''1.375'
$ws.Range("E29").Value = '  +1.20%  '
# Row 30: Hedera
$ws.Range("D30").Value = '''0.05924'
$ws.Range("E30").Value = '  -4.78%  '
# Row 31: PancakeSwap
$ws.Range("E31").Value = '  -3.05%  '
# Row 32: Filecoin
$ws.Range("D32").Value = '''3.332'
$ws.Range("E32").Value = '  -3.36%  '
# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").Value = '''3.325'
$ws.Range("E33").Value = '  -4.84%  '
# Row 34: LidoDAOToken
$ws.Range("D34").Value = '''1.598'
$ws.Range("E34").Value = '  -2.00%  '
# Row 35: ARBITRUM
$ws.Range("D35").Value = '''0.9648'
$ws.Range("E35").Value = '  -3.83%  '
# Row 36: HuobiToken
$ws.Range("D36").Value = '''2.381'
$ws.Range("E36").Value = '  -0.76%  '
# Row 37: MXToken
$ws.Range("D37").Value = '''2.735'
$ws.Range("E37").Value = '  +0.75%  '
# Row 38: ImmutableX
$ws.Range("D38").Value = '''0.5814'
$ws.Range("E38").Value = '  -4.12%  '
# Row 39: VeChain
$ws.Range("D39").Value = '''0.01592'
$ws.Range("E39").Value = '  -1.57%  '
# Row 40: PaxDollar
$ws.Range("E40").Value = '  -0.31%  '
# Row 41: TrustWalletToken
$ws.Range("D41").Value = '''0.8383'
$ws.Range("E41").Value = '  -3.38%  '
# Row 42: Maker
$ws.Range("D42").Value = '1.035.36'
$ws.Range("E42").Value = '  -4.21%  '
# Row 43: FraxShare
$ws.Range("D43").Value = '''5.672'
$ws.Range("E43").Value = '  -7.31%  '
# Row 44: Quant
$ws.Range("E44").Value = '  -0.68%  '
# Row 45: RocketPoolETH
$ws.Range("D45").Value = '1.788.76'
$ws.Range("E45").Value = '  -1.95%  '
# Row 46: BabyDogeCoin
$ws.Range("E46").Value = '  -2.21%  '
# Row 47: Frax
$ws.Range("D47").Value = '''0.9977'
$ws.Range("E47").Value = '  -0.36%  '
# Row 48: Aave
$ws.Range("D48").Value = '''54.22'
$ws.Range("E48").Value = '  -3.62%  '
# Row 49: EnergySwap
$ws.Range("D49").Value = '''7.928'
$ws.Range("E49").Value = '  -1.29%  '
# Row 50: Cronos
$ws.Range("E50").Value = '  -1.06%  '
# Row 51: Mantle
$ws.Range("D51").Value = '''0.4230'
$ws.Range("E51").Value = '  -0.56%  '
